$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.295.60"
$ws.Range("E2").Value = "  +3.48%  "
$ws.Range("D3").Value = "3.062.53"
$ws.Range("E3").Value = "  +5.15%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.38%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +3.37%  "
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("E10").Value = "  +4.25%  "
$ws.Range("E11").Value = "  +6.86%  "
$ws.Range("D12").Value = "3.589.74"
$ws.Range("E12").Value = "  +5.35%  "
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("E15").Value = "  +4.11%  "
$ws.Range("D16").Value = "57.340.31"
$ws.Range("E16").Value = "  +3.56%  "
$ws.Range("D17").Value = "3.066.59"
$ws.Range("E17").Value = "  +5.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "337.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.65%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("E23").Value = "  +3.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.68%  "
$ws.Range("E25").Value = "  +7.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "0.0₃0933"
$ws.Range("E27").Value = "  +12.44%  "
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("E31").Value = "  +4.91%  "
$ws.Range("E32").Value = "  +5.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("E35").Value = "  +5.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "25.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.57%  "
$ws.Range("E38").Value = "  +4.86%  "
$ws.Range("D39").Value = "3.101.51"
$ws.Range("E39").Value = "  +5.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("E41").Value = "  +4.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.668"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.45%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "2.248.69"
$ws.Range("E44").Value = "  +7.22%  "
$ws.Range("E45").Value = "  +8.77%  "
$ws.Range("E46").Value = "  +4.73%  "
$ws.Range("E47").Value = "  +4.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.28%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0865"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.68%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.15%  "
